$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The password text itself was corrected (jaga@298 -> jaga@2983). D2 and the
# new D3 cell both point at this same shared string, so updating D2 in place
# (rather than leaving it alone) reproduces the shared-string edit.
$ws.Range("D2").Value = "jaga@2983"

# Add a new row 3 that mirrors row 2 (TestCase/UserName values + SignIn hyperlinks).
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2

# C3: same hyperlink target/display text as C2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:testjaga002@gmail.com")
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "testjaga002@gmail.com"

# D3: hyperlink to the (corrected) password, formatted like D2.
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:jaga@2983")
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value = "jaga@2983"

$null = $ws.Range("D2").Select()
